$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 211, shifting existing rows 211-275 down to 212-276.
$ws.Rows(211).Insert()

# Populate the newly inserted row 211 with the new record.
$ws.Range("A211").Value = 3
$ws.Range("B211").Value = "Femacal de La Calera"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 44468
$ws.Range("E211").Value = 5
$ws.Range("F211").Value = 100112021
$ws.Range("G211").Value = "Ají"
$ws.Range("H211").Value = "Americana (o)"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 60
$ws.Range("K211").Value = 41000
$ws.Range("L211").Value = 42000
$ws.Range("M211").Value = 41417
$ws.Range("N211").Value = "`$/caja 15 kilos"
$ws.Range("O211").Value = "Región de Arica y Parinacota"
$ws.Range("P211").Value = 2761
$ws.Range("Q211").Value = 15
$ws.Range("R211").Value = "Hortaliza"
